$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out row 2 extra data and all of row 3, keep only C2 = 10000 (no style)
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("A2:G2").ClearFormats()

$ws.Range("C2").Value = 10000

$ws.Range("A3:G3").ClearContents()
$ws.Range("A3:G3").ClearFormats()

# Remove the date-format style from G1 and column G
$ws.Range("G1").ClearFormats()
$ws.Range("G1").Value = "Holding"
$ws.Columns.Item(7).ClearFormats()

# Adjust column widths
$ws.Columns.Item(1).ColumnWidth = 5.140625

# Adjust the window size recorded in the workbook view
$excel.ActiveWindow.Width = 15345
$excel.ActiveWindow.Height = 6705
